# Practical creation of a table (MYSQL "friends" table) added to the
# "Commands" sheet of the MySQL and Mongodb Queries workbook.
#
# Semantic changes (per the authoritative diff):
#   A5  : "Establishing connection to mysql database" -> "Establishing connection"
#   B10 : "...mycursor.execute(\"CREATE TABLE customers (name VARCHAR(255), address VARCHAR(255))\")"
#         -> richer text referencing a new "friends" table, with the CREATE TABLE
#            statement rendered in bold.
#   C10 : style only changes (Consolas/no-border style -> the regular
#         bordered/wrapped body style used elsewhere in the sheet); text unchanged.
#   D10 : was empty -> new note about data types / display width, with "MYSQL:"
#         bold+underlined and "Note from stackoverflow" bold.
#   Row 10 height grows from 60 to 180 to fit the new D10 content.
#   Frozen pane / active selection move from A5/C7 down to A10/C10.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Commands")

# --- A5: shorten the "establishing connection" label ------------------------
$ws.Range("A5").Value = "Establishing connection"

# --- B10: new CREATE TABLE example (friends table) --------------------------
$b10Run1 = "Mention the database name while establishing connection to mysql databale schema`nmycursor.execute(`n`""
$b10Run2 = "CREATE TABLE friends (id INT AUTO_INCREMENT PRIMARY KEY, name VARCHAR(255), email VARCHAR(255), mobile INT)"
$b10Run3 = "`"`n)"
$b10Full = $b10Run1 + $b10Run2 + $b10Run3

$ws.Range("B10").Value = $b10Full

$b10Start1 = 1
$b10Start2 = $b10Start1 + $b10Run1.Length
$b10Start3 = $b10Start2 + $b10Run2.Length

$b10Chars2 = $ws.Range("B10").Characters($b10Start2, $b10Run2.Length)
$b10Chars2.Font.Bold = $true

# --- C10: keep the text, only normalise the cell style -----------------------
# Copy the plain wrapped/bordered body style (used by e.g. C9) onto C10, then
# drop the special "code" (Consolas) styling it used to carry.
$ws.Range("C9").Copy() | Out-Null
$ws.Range("C10").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# --- D10: brand-new note about data types / display width -------------------
$d10Run1 = "MYSQL:"
$d10Run2 = "`n1) Create table - For data types and widh details please refer online. Basically 3 datatypes.`n`"`"`""
$d10Run3 = "`nNote from stackoverflow"
$d10Run4 = ": Display width doesn't change storage requirements for a data type.`nDisplay width doesn't alter the actual data in any way (ie: it stores the entire value for the data)`nA column returns it's full value when called in a query, regardless of the display width (the book directly contradicts this claim it makes as seen above)`n`"`"`""
$d10Full = $d10Run1 + $d10Run2 + $d10Run3 + $d10Run4

$ws.Range("D10").Value = $d10Full

$d10Start1 = 1
$d10Start2 = $d10Start1 + $d10Run1.Length
$d10Start3 = $d10Start2 + $d10Run2.Length
$d10Start4 = $d10Start3 + $d10Run3.Length

$d10Chars1 = $ws.Range("D10").Characters($d10Start1, $d10Run1.Length)
$d10Chars1.Font.Bold = $true
$d10Chars1.Font.Underline = $true

$d10Chars3 = $ws.Range("D10").Characters($d10Start3, $d10Run3.Length)
$d10Chars3.Font.Bold = $true

# --- Row 10 height grows to fit the new note text ----------------------------
$ws.Range("A10:E10").RowHeight = 180

# --- Frozen pane / selection follow the edited rows --------------------------
$ws.Range("C10").Select()
$excel.ActiveWindow.FreezePanes = $false
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("C10").Select()
